# Updated symbol list on Tue Dec 20 19:40:23 UTC 2022 with GitHub Actions
#
# Every "Price" (column D) cell in the sheet is stored as TEXT (the scraper
# wrote it that way), so an update like "250.14" -> "250.05" must stay text,
# not become the number 250.05 (which would lose e.g. trailing zeros such as
# "0.04050" or the long decimals like "0.00000000749"). Flip the cell to the
# Text number format before assigning the new literal, then clear the
# number-format back off the cell so no visible formatting change lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.ClearFormats()
}

Set-TextValue "D2"  "250.05"
Set-TextValue "D3"  "22.83"
Set-TextValue "D4"  "5.443"
Set-TextValue "D5"  "0.05653"
Set-TextValue "D8"  "0.8155"
Set-TextValue "D9"  "0.9286"
Set-TextValue "D11" "0.07551"
Set-TextValue "D12" "0.03166"
Set-TextValue "D13" "0.03108"
Set-TextValue "D14" "0.09355"
Set-TextValue "D15" "3.559"
Set-TextValue "D16" "0.001588"
Set-TextValue "D18" "0.0005783"
Set-TextValue "E18" "17OneONE"
Set-TextValue "D19" "0.006384"
Set-TextValue "D20" "0.004999"
Set-TextValue "D21" "0.001029"
Set-TextValue "D22" "0.0001498"
Set-TextValue "D24" "2.202"
Set-TextValue "D25" "0.3299"
Set-TextValue "E27" "26AAXTokenAABWorstin24h"
Set-TextValue "D28" "0.0003028"
Set-TextValue "D40" "0.04048"
Set-TextValue "D41" "0.006771"
Set-TextValue "D42" "0.1070"
Set-TextValue "D43" "0.002706"
Set-TextValue "D44" "0.007551"
Set-TextValue "D45" "0.00005795"
Set-TextValue "D46" "0.00000000749"
Set-TextValue "D47" "0.4994"
Set-TextValue "D48" "0.2413"
Set-TextValue "D49" "0.00002098"
